# Updated capital structure database
# - Malawi "Bank (Money Center)" peer set refreshed:
#     * row 2 (unnamed company "3") relabeled "2" and all metrics refreshed
#     * row 3 "NBS Bank Plc (MAL:NBS)" replaced by "National Bank of Malawi plc (MAL:NBM)" with its metrics
#     * row 4 "Standard Bank Limited (MAL:STANDARD)" replaced by "FDH Bank Plc (MAL:FDHB)" with its metrics
#     * old row 5 ("National Bank of Malawi plc (MAL:NBM)") is removed, shrinking the table to 4 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 5 entirely; this also shrinks the sheet dimension to A1:AQ4
$ws.Rows("5:5").Delete()

# --- Row 2 -------------------------------------------------------------
$ws.Range("B2").Value = "'2"
$ws.Range("D2").Value = 0.159
$ws.Range("E2").Value = 0.11
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 33.6
$ws.Range("L2").Value = 0.2459736456808199
$ws.Range("M2").Value = 6.98
$ws.Range("N2").Value = 0.01338960291578746
$ws.Range("O2").Value = 0.2077380952380952
$ws.Range("P2").Value = 6.98
$ws.Range("Q2").Value = 0.01338960291578746
$ws.Range("R2").Value = 0.2077380952380952
$ws.Range("U2").Value = 51.40000000000001
$ws.Range("V2").Value = 0.09859965470938041
$ws.Range("W2").Value = 0.3110145568328685
$ws.Range("X2").Value = 0.0797775747302297
$ws.Range("Y2").Value = 0.2312369821026388
$ws.Range("Z2").Value = 1.137954015328224
$ws.Range("AB2").Value = 0.07923966412932001
$ws.Range("AC2").Value = -0.07923966412932001
$ws.Range("AD2").Value = 51.1
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 51.1
$ws.Range("AG2").Value = -0.3000000000000043
$ws.Range("AH2").Value = 0.08927323549965061
$ws.Range("AI2").Value = 0.2193133047210301
$ws.Range("AJ2").Value = -0.0005758157389635398
$ws.Range("AK2").Value = -0.001651982378854649
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3 -------------------------------------------------------------
$ws.Range("B3").Value = "National Bank of Malawi plc (MAL:NBM)"
$ws.Range("K3").Value = 23.1
$ws.Range("L3").Value = 0.257238307349666
$ws.Range("M3").Value = 5.57
$ws.Range("N3").Value = 0.01419469928644241
$ws.Range("O3").Value = 0.2411255411255411
$ws.Range("P3").Value = 5.57
$ws.Range("Q3").Value = 0.01419469928644241
$ws.Range("R3").Value = 0.2411255411255411
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 35.6
$ws.Range("V3").Value = 0.09072375127420999
$ws.Range("W3").Value = 0.191701244813278
$ws.Range("X3").Value = 0.07653263110929406
$ws.Range("Y3").Value = 0.115168613703984
$ws.Range("Z3").Value = 0.9852973447443493
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.07654141338629493
$ws.Range("AC3").Value = -0.07654141338629493
$ws.Range("AD3").Value = 24.5
$ws.Range("AF3").Value = 24.5
$ws.Range("AG3").Value = -11.1
$ws.Range("AH3").Value = 0.05876709042935956
$ws.Range("AI3").Value = 0.14
$ws.Range("AJ3").Value = -0.02911093627065304
$ws.Range("AK3").Value = -0.07962697274031565

# --- Row 4 -------------------------------------------------------------
$ws.Range("B4").Value = "FDH Bank Plc (MAL:FDHB)"
$ws.Range("D4").Value = 0.159
$ws.Range("E4").Value = 0.11
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 10.5
$ws.Range("L4").Value = 0.2243589743589744
$ws.Range("M4").Value = 1.41
$ws.Range("N4").Value = 0.01093871217998448
$ws.Range("O4").Value = 0.1342857142857143
$ws.Range("P4").Value = 1.41
$ws.Range("Q4").Value = 0.01093871217998448
$ws.Range("R4").Value = 0.1342857142857143
$ws.Range("U4").Value = 15.8
$ws.Range("V4").Value = 0.1225756400310318
$ws.Range("W4").Value = 0.4303278688524591
$ws.Range("X4").Value = 0.08302251835116532
$ws.Range("Y4").Value = 0.3473053505012937
$ws.Range("Z4").Value = 1.619377162629758
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.08193791487234511
$ws.Range("AC4").Value = -0.08193791487234511
$ws.Range("AD4").Value = 26.6
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 26.6
$ws.Range("AG4").Value = 10.8
$ws.Range("AH4").Value = 0.1710610932475884
$ws.Range("AI4").Value = 0.4586206896551724
$ws.Range("AJ4").Value = 0.07730851825340014
$ws.Range("AK4").Value = 0.2559241706161137
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
